# Auto-generated edit script: apply numeric updates to the "Leve Profit" data
# sheets (ALC, ARM, CRP, CUL, GSM, LTW, WVR) per the scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

# ALC row 64
$ws = $wb.Worksheets.Item(1)
$ws.Range("H64").Value = 4731.6665
$ws.Range("I64").Value = 4597.5
$ws.Range("J64").Value = 5000
$ws.Range("K64").Value = 4597.5
$ws.Range("L64").Value = 5000
$ws.Range("M64").Value = -4349.5
$ws.Range("N64").Value = -5496

# ALC row 67
$ws = $wb.Worksheets.Item(1)
$ws.Range("H67").Value = 4731.6665
$ws.Range("I67").Value = 4597.5
$ws.Range("J67").Value = 5000
$ws.Range("K67").Value = 4597.5
$ws.Range("L67").Value = 5000
$ws.Range("M67").Value = -3739.5
$ws.Range("N67").Value = -6716

# ALC row 74
$ws = $wb.Worksheets.Item(1)
$ws.Range("H74").Value = 8632.666999999999
$ws.Range("I74").Value = 8632.666999999999
$ws.Range("K74").Value = 8632.666999999999
$ws.Range("M74").Value = -7696.666999999999

# ALC row 77
$ws = $wb.Worksheets.Item(1)
$ws.Range("H77").Value = 8632.666999999999
$ws.Range("I77").Value = 8632.666999999999
$ws.Range("K77").Value = 43163.335
$ws.Range("M77").Value = -38483.335

# ALC row 101
$ws = $wb.Worksheets.Item(1)
$ws.Range("H101").Value = 5670
$ws.Range("I101").Value = 3962.5
$ws.Range("K101").Value = 11887.5
$ws.Range("M101").Value = -10265.5

# ALC row 132
$ws = $wb.Worksheets.Item(1)
$ws.Range("H132").Value = 2828.6553
$ws.Range("I132").Value = 2462.7693
$ws.Range("J132").Value = 5999.6665
$ws.Range("K132").Value = 7388.3079
$ws.Range("L132").Value = 17998.9995
$ws.Range("M132").Value = -4858.3079
$ws.Range("N132").Value = -23058.9995

# ALC row 138
$ws = $wb.Worksheets.Item(1)
$ws.Range("H138").Value = 2810.8289
$ws.Range("I138").Value = 834.375
$ws.Range("J138").Value = 3337.8833
$ws.Range("K138").Value = 2503.125
$ws.Range("L138").Value = 10013.6499
$ws.Range("M138").Value = 2636.875
$ws.Range("N138").Value = -20293.6499

# ARM row 7
$ws = $wb.Worksheets.Item(2)
$ws.Range("H7").Value = 151000
$ws.Range("J7").Value = 151000
$ws.Range("L7").Value = 151000
$ws.Range("N7").Value = -151228

# ARM row 27
$ws = $wb.Worksheets.Item(2)
$ws.Range("H27").Value = 7895
$ws.Range("J27").Value = 7895
$ws.Range("L27").Value = 7895
$ws.Range("N27").Value = -8263

# ARM row 32
$ws = $wb.Worksheets.Item(2)
$ws.Range("H32").Value = 17864692
$ws.Range("I32").Value = 17864692
$ws.Range("K32").Value = 17864692
$ws.Range("M32").Value = -17864405

# ARM row 38
$ws = $wb.Worksheets.Item(2)
$ws.Range("H38").Value = 5019
$ws.Range("I38").Value = 5019
$ws.Range("K38").Value = 5019
$ws.Range("M38").Value = -4552

# ARM row 61
$ws = $wb.Worksheets.Item(2)
$ws.Range("H61").Value = 20881484
$ws.Range("I61").Value = 45459030
$ws.Range("J61").Value = 85096.84
$ws.Range("K61").Value = 45459030
$ws.Range("L61").Value = 85096.84
$ws.Range("M61").Value = -45458818
$ws.Range("N61").Value = -85520.84

# ARM row 74
$ws = $wb.Worksheets.Item(2)
$ws.Range("H74").Value = 13899130
$ws.Range("I74").Value = 22728318
$ws.Range("J74").Value = 24692.143
$ws.Range("K74").Value = 22728318
$ws.Range("L74").Value = 24692.143
$ws.Range("M74").Value = -22727444
$ws.Range("N74").Value = -26440.143

# ARM row 77
$ws = $wb.Worksheets.Item(2)
$ws.Range("H77").Value = 13899130
$ws.Range("I77").Value = 22728318
$ws.Range("J77").Value = 24692.143
$ws.Range("K77").Value = 113641590
$ws.Range("L77").Value = 123460.715
$ws.Range("M77").Value = -113637222
$ws.Range("N77").Value = -132196.715

# ARM row 110
$ws = $wb.Worksheets.Item(2)
$ws.Range("H110").Value = 11000
$ws.Range("I110").Value = 2000
$ws.Range("J110").Value = 20000
$ws.Range("K110").Value = 2000
$ws.Range("L110").Value = 20000
$ws.Range("M110").Value = 45
$ws.Range("N110").Value = -24090

# ARM row 132
$ws = $wb.Worksheets.Item(2)
$ws.Range("H132").Value = 5799.5864
$ws.Range("I132").Value = 3119.682
$ws.Range("K132").Value = 9359.045999999998
$ws.Range("M132").Value = -6829.045999999998

# ARM row 136
$ws = $wb.Worksheets.Item(2)
$ws.Range("H136").Value = 20881484
$ws.Range("I136").Value = 45459030
$ws.Range("J136").Value = 85096.84
$ws.Range("K136").Value = 136377090
$ws.Range("L136").Value = 255290.52
$ws.Range("M136").Value = -136374540
$ws.Range("N136").Value = -260390.52

# CRP row 7
$ws = $wb.Worksheets.Item(4)
$ws.Range("H7").Value = 2793.4614
$ws.Range("I7").Value = 176.77777
$ws.Range("K7").Value = 176.77777
$ws.Range("M7").Value = -63.77777

# CRP row 31
$ws = $wb.Worksheets.Item(4)
$ws.Range("H31").Value = 730858.9
$ws.Range("I31").Value = 1655.2
$ws.Range("K31").Value = 1655.2
$ws.Range("M31").Value = -1360.2

# CRP row 32
$ws = $wb.Worksheets.Item(4)
$ws.Range("H32").Value = 2010
$ws.Range("I32").Value = 2010
$ws.Range("K32").Value = 2010
$ws.Range("M32").Value = -1694

# CRP row 34
$ws = $wb.Worksheets.Item(4)
$ws.Range("H34").Value = 730858.9
$ws.Range("I34").Value = 1655.2
$ws.Range("K34").Value = 1655.2
$ws.Range("M34").Value = -1453.2

# CRP row 58
$ws = $wb.Worksheets.Item(4)
$ws.Range("H58").Value = 1872.5714
$ws.Range("I58").Value = 1718.6666
$ws.Range("J58").Value = 2796
$ws.Range("K58").Value = 1718.6666
$ws.Range("L58").Value = 2796
$ws.Range("M58").Value = -1515.6666
$ws.Range("N58").Value = -3202

# CRP row 99
$ws = $wb.Worksheets.Item(4)
$ws.Range("H99").Value = 3755.75
$ws.Range("I99").Value = 3007.6667
$ws.Range("J99").Value = 6000
$ws.Range("K99").Value = 3007.6667
$ws.Range("L99").Value = 6000
$ws.Range("M99").Value = -1509.6667
$ws.Range("N99").Value = -8996

# CRP row 126
$ws = $wb.Worksheets.Item(4)
$ws.Range("H126").Value = 3755.75
$ws.Range("I126").Value = 3007.6667
$ws.Range("J126").Value = 6000
$ws.Range("K126").Value = 9023.000100000001
$ws.Range("L126").Value = 18000
$ws.Range("M126").Value = -6553.000100000001
$ws.Range("N126").Value = -22940

# CRP row 136
$ws = $wb.Worksheets.Item(4)
$ws.Range("H136").Value = 1872.5714
$ws.Range("I136").Value = 1718.6666
$ws.Range("J136").Value = 2796
$ws.Range("K136").Value = 5155.9998
$ws.Range("L136").Value = 8388
$ws.Range("M136").Value = -2605.9998
$ws.Range("N136").Value = -13488

# CUL row 11
$ws = $wb.Worksheets.Item(5)
$ws.Range("H11").Value = 353.7857
$ws.Range("I11").Value = 273.30435
$ws.Range("K11").Value = 819.91305
$ws.Range("M11").Value = -679.91305

# CUL row 18
$ws = $wb.Worksheets.Item(5)
$ws.Range("H18").Value = 53
$ws.Range("I18").Value = 53
$ws.Range("K18").Value = 159
$ws.Range("M18").Value = 10

# CUL row 140
$ws = $wb.Worksheets.Item(5)
$ws.Range("H140").Value = 274539.1
$ws.Range("I140").Value = 274539.1
$ws.Range("K140").Value = 823617.2999999999
$ws.Range("M140").Value = -818437.2999999999

# GSM row 5
$ws = $wb.Worksheets.Item(6)
$ws.Range("H5").Value = 5837.75
$ws.Range("I5").Value = 3560.2666
$ws.Range("J5").Value = 40000
$ws.Range("K5").Value = 3560.2666
$ws.Range("L5").Value = 40000
$ws.Range("M5").Value = -3448.2666
$ws.Range("N5").Value = -40224

# GSM row 62
$ws = $wb.Worksheets.Item(6)
$ws.Range("H62").Value = 99708.25
$ws.Range("J62").Value = 99708.25
$ws.Range("L62").Value = 99708.25
$ws.Range("N62").Value = -101080.25

# GSM row 63
$ws = $wb.Worksheets.Item(6)
$ws.Range("H63").Value = 49995
$ws.Range("J63").Value = 49995
$ws.Range("L63").Value = 49995
$ws.Range("N63").Value = -51367

# GSM row 65
$ws = $wb.Worksheets.Item(6)
$ws.Range("H65").Value = 99708.25
$ws.Range("J65").Value = 99708.25
$ws.Range("L65").Value = 299124.75
$ws.Range("N65").Value = -305988.75

# GSM row 66
$ws = $wb.Worksheets.Item(6)
$ws.Range("H66").Value = 49995
$ws.Range("J66").Value = 49995
$ws.Range("L66").Value = 149985
$ws.Range("N66").Value = -156849

# GSM row 126
$ws = $wb.Worksheets.Item(6)
$ws.Range("H126").Value = 2999.3333
$ws.Range("I126").Value = 2999.5
$ws.Range("J126").Value = 2999
$ws.Range("K126").Value = 8998.5
$ws.Range("L126").Value = 8997
$ws.Range("M126").Value = -6528.5
$ws.Range("N126").Value = -13937

# GSM row 132
$ws = $wb.Worksheets.Item(6)
$ws.Range("H132").Value = 71431520
$ws.Range("I132").Value = 71431520
$ws.Range("K132").Value = 214294560
$ws.Range("M132").Value = -214292030

# LTW row 17
$ws = $wb.Worksheets.Item(7)
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("N17").ClearContents()

# LTW row 22
$ws = $wb.Worksheets.Item(7)
$ws.Range("H22").Value = 3024.75
$ws.Range("I22").Value = 3766.6667
$ws.Range("K22").Value = 3766.6667
$ws.Range("M22").Value = -3471.6667

# LTW row 27
$ws = $wb.Worksheets.Item(7)
$ws.Range("H27").Value = 3024.75
$ws.Range("I27").Value = 3766.6667
$ws.Range("K27").Value = 3766.6667
$ws.Range("M27").Value = -3659.6667

# LTW row 46
$ws = $wb.Worksheets.Item(7)
$ws.Range("H46").Value = 2570.2964
$ws.Range("J46").Value = 3256.125
$ws.Range("L46").Value = 3256.125
$ws.Range("N46").Value = -3632.125

# LTW row 122
$ws = $wb.Worksheets.Item(7)
$ws.Range("H122").Value = 6430.9375
$ws.Range("I122").Value = 5710.5557
$ws.Range("K122").Value = 17131.6671
$ws.Range("M122").Value = -14681.6671

# LTW row 136
$ws = $wb.Worksheets.Item(7)
$ws.Range("H136").Value = 67430.336
$ws.Range("I136").Value = 40607.58
$ws.Range("K136").Value = 121822.74
$ws.Range("M136").Value = -119272.74

# WVR row 113
$ws = $wb.Worksheets.Item(8)
$ws.Range("H113").Value = 812
$ws.Range("I113").Value = 778.2857
$ws.Range("J113").Value = 930
$ws.Range("K113").Value = 2334.8571
$ws.Range("L113").Value = 2790
$ws.Range("M113").Value = -164.8571000000002
$ws.Range("N113").Value = -7130

# WVR row 122
$ws = $wb.Worksheets.Item(8)
$ws.Range("H122").Value = 4754.615
$ws.Range("I122").Value = 3579.3845
$ws.Range("J122").Value = 5929.846
$ws.Range("K122").Value = 10738.1535
$ws.Range("L122").Value = 17789.538
$ws.Range("M122").Value = -8288.1535
$ws.Range("N122").Value = -22689.538
